# Edit script for "Modelo CFD DroneBoX/Datos Aero DBX v2.xlsx"
# Raw Data sheet (sheet4.xml): insert one new data row (alfa_deg = 90,
# beta_deg = 0 ...) right before the existing "alfa=0,beta=0..." block
# (old row 23), which pushes the rest of the table down by one row, and
# append one brand-new row at the very end (alfa_deg = -90).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Raw Data")

# --- 1. Insert a new row before the current row 23 (shifts 23..129 -> 24..130) ---
$ws.Rows.Item(23).Insert()

# --- 2. Populate the newly inserted row 23 ---
$ws.Range("A23").Value = 90
$ws.Range("B23").Value = 0
$ws.Range("C23").Value = 0
$ws.Range("D23").Value = 0
$ws.Range("E23").Value = 0
$ws.Range("F23").Value = [double]"0.40509845602615019"
$ws.Range("G23").Value = [double]"-8.3563578592191351E-3"
$ws.Range("H23").Value = [double]"-1.7614898218215362"
$ws.Range("I23").Value = [double]"-2.4137552639073067E-3"
$ws.Range("J23").Value = [double]"-1.2480834362936954"
$ws.Range("K23").Value = [double]"1.1477498080717985E-3"
$ws.Range("L23").Value = [double]"0.4050984560261503"
$ws.Range("M23").Value = [double]"1.7614898218215362"

# --- 3. Append a brand-new row 131 at the bottom of the table ---
$ws.Range("A131").Value = -90
$ws.Range("B131").Value = 0
$ws.Range("C131").Value = 0
$ws.Range("D131").Value = 0
$ws.Range("E131").Value = 0
$ws.Range("F131").Value = [double]"0.15316996570630381"
$ws.Range("G131").Value = [double]"1.1236940653992767E-2"
$ws.Range("H131").Value = [double]"1.8409158950242044"
$ws.Range("I131").Value = [double]"9.0540662991842619E-3"
$ws.Range("J131").Value = [double]"1.0436648502996442"
$ws.Range("K131").Value = [double]"-6.4741679085321576E-4"
$ws.Range("L131").Value = [double]"-0.15316996570630392"
$ws.Range("M131").Value = [double]"1.8409158950242044"

# --- 4. Fix up the sheet view selection to match the new used range ---
$ws.Range("K1:M131").Select()
